# Daily attendance processing - 2025-10-30 06:55:06
#
# Per-session-row corrections to the attendance/session-analysis sheet:
#   - "Recorded By" (col G): re-order the recorder list so the backup
#     service account sorts first (e.g. "System, backup@backdoor.com"
#     -> "backup@backdoor.com, System").
#   - "Students" (col H) and the class-statistics roll-up (cols L, M, S):
#     the class roster size changed for a couple of groups, so every
#     "x/old_total" attendance fraction and the dependent summary cells
#     (roster size, recorded-rate %) are refreshed to "x/new_total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write a literal text value (e.g. a percentage string like "71.2%")
# into a cell without Excel's number-parsing reinterpreting it as a
# numeric percentage. Goes through a throwaway formula (a quoted-string
# literal always evaluates to text) then flattens the cell back to a
# plain value in place via copy / paste-special-values, which drops the
# formula but keeps the cell's existing style untouched.
function Set-LiteralText($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.Formula = '="' + $text + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)  # xlPasteValues
}

$ws.Range("G2").Value = 'backup@backdoor.com, system, System'
$ws.Range("H2").Value = '33/52'
$ws.Range("H3").Value = '52/52'
$ws.Range("G4").Value = 'backup@backdoor.com, System'
$ws.Range("H4").Value = '42/52'
$ws.Range("G5").Value = 'backup@backdoor.com, System'
$ws.Range("H5").Value = '35/52'
$ws.Range("H6").Value = '45/52'
$ws.Range("H7").Value = '35/52'
$ws.Range("G8").Value = 'backup@backdoor.com, System'
$ws.Range("H8").Value = '44/52'
$ws.Range("H9").Value = '21/52'
$ws.Range("H10").Value = '27/52'
Set-LiteralText "L10" '71.2%'
$ws.Range("H11").Value = '34/52'
$ws.Range("H12").Value = '31/52'
$ws.Range("H13").Value = '36/52'
$ws.Range("H14").Value = '36/52'
$ws.Range("H15").Value = '43/52'
$ws.Range("M15").Value = 52
Set-LiteralText "S15" '70.0%'
$ws.Range("H16").Value = '42/52'
$ws.Range("H17").Value = '36/52'
$ws.Range("H18").Value = '38/52'
$ws.Range("H19").Value = '35/52'
$ws.Range("M19").Value = 56
Set-LiteralText "S19" '73.0%'
$ws.Range("H20").Value = '42/52'
$ws.Range("H21").Value = '42/52'
$ws.Range("H22").Value = '31/52'
$ws.Range("H23").Value = '41/52'
$ws.Range("H24").Value = '27/52'
$ws.Range("H25").Value = '26/52'
$ws.Range("H26").Value = '0/52'
$ws.Range("H27").Value = '0/52'
$ws.Range("H28").Value = '0/52'
$ws.Range("G29").Value = 'backup@backdoor.com, system, System'
$ws.Range("G31").Value = 'backup@backdoor.com, System'
$ws.Range("G32").Value = 'backup@backdoor.com, System'
$ws.Range("G35").Value = 'backup@backdoor.com, System'
$ws.Range("G56").Value = 'backup@backdoor.com, system, System'
$ws.Range("G58").Value = 'backup@backdoor.com, System'
$ws.Range("G59").Value = 'backup@backdoor.com, System'
$ws.Range("G62").Value = 'backup@backdoor.com, System'
$ws.Range("G83").Value = 'backup@backdoor.com, System'
$ws.Range("G84").Value = 'backup@backdoor.com, System'
$ws.Range("G85").Value = 'backup@backdoor.com, System'
$ws.Range("G109").Value = 'backup@backdoor.com, System'
$ws.Range("H109").Value = '30/56'
$ws.Range("G110").Value = 'backup@backdoor.com, System'
$ws.Range("H110").Value = '38/56'
$ws.Range("G111").Value = 'backup@backdoor.com, System'
$ws.Range("H111").Value = '55/56'
$ws.Range("H112").Value = '54/56'
$ws.Range("H113").Value = '26/56'
$ws.Range("H114").Value = '37/56'
$ws.Range("H115").Value = '39/56'
$ws.Range("H116").Value = '54/56'
$ws.Range("H117").Value = '54/56'
$ws.Range("H118").Value = '47/56'
$ws.Range("H119").Value = '42/56'
$ws.Range("H120").Value = '45/56'
$ws.Range("H121").Value = '43/56'
$ws.Range("H122").Value = '41/56'
$ws.Range("H123").Value = '36/56'
$ws.Range("H124").Value = '38/56'
$ws.Range("H125").Value = '45/56'
$ws.Range("H126").Value = '36/56'
$ws.Range("H127").Value = '33/56'
$ws.Range("H128").Value = '42/56'
$ws.Range("H129").Value = '29/56'
$ws.Range("H130").Value = '35/56'
$ws.Range("H131").Value = '0/56'
$ws.Range("H132").Value = '0/56'
$ws.Range("H133").Value = '0/56'
$ws.Range("H134").Value = '0/56'
$ws.Range("G135").Value = 'backup@backdoor.com, System'
$ws.Range("G136").Value = 'backup@backdoor.com, System'
$ws.Range("G137").Value = 'backup@backdoor.com, System'
